# Apply imputed / cleared values to the missing_data worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Individual cell corrections (rows keep their original numbering here,
#     row deletions for RM 232 / SC 92 are performed afterwards) ---

$ws.Range("C2").Value = 14.9
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E6").Value = -5.7
$ws.Range("D8").Value = -13.9
$ws.Range("D10").Value = -14.7
$ws.Range("C11").Value = 11.4
$ws.Range("E11").Value = -7.9
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("E13").Value = -5.3
$ws.Range("D15").Value = -15.2
$ws.Range("E17").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("E18").Value = -8.5
$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("C21").Value = 12.7
$ws.Range("E24").ClearContents()
$ws.Range("C25").ClearContents()
$ws.Range("D25").Value = -15.5
$ws.Range("E25").Value = -7.1

# Rows below the two rows that will be removed (RM 232 at row 26,
# SC 92 at row 28) - edited here while row numbers still match the
# pre-deletion layout.
$ws.Range("D29").Value = -14.6
$ws.Range("B31").ClearContents()
$ws.Range("D31").ClearContents()
$ws.Range("E33").ClearContents()
$ws.Range("E34").ClearContents()
$ws.Range("B35").Value = -19.5
$ws.Range("C35").Value = 10.4
$ws.Range("D35").ClearContents()

# --- Remove the two rows entirely (RM 232 and SC 92) ---
# Delete the lower-numbered row last so the higher row index is still valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()
